# Refresh the "cryptos" price table (GitHub Actions scheduled update).
# Only column D (Price) and column E (Volume(1h)) move; everything else is
# left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free inline updates below. Price strings such as "1.000" or "0.7755"
# look numeric to Excel's automatic type conversion, so a bare
# `.Value = "1.000"` would silently store the Double 1 instead of the text
# "1.000" (losing the trailing zero(s) the source feed renders). We force
# text entry the same way typing an apostrophe-prefixed value in the UI
# does, then reset `.Style` back to the workbook's built-in Normal style so
# the cell does not end up tagged with an explicit "quote prefix" style that
# was never part of the original formatting.

$ws.Range("D2").Value = '''29.933.02'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.01%  '
$ws.Range("D3").Value = '''1.894.92'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.24%  '
$ws.Range("D4").Value = '''1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''0.7770'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.37%  '
$ws.Range("D6").Value = '''244.66'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").Value = '''1.001'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '''0.3139'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.31%  '
$ws.Range("D9").Value = '''25.86'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.09%  '
$ws.Range("D10").Value = '''0.07251'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.68%  '
$ws.Range("D11").Value = '''0.09009'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +10.92%  '
$ws.Range("D12").Value = '''0.7729'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = '''1.947.49'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.18%  '
$ws.Range("D14").Value = '''5.476'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.43%  '
$ws.Range("D15").Value = '''94.84'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.08%  '
$ws.Range("D16").Value = '''6.217'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.22%  '
$ws.Range("D17").Value = '''30.050.04'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.51%  '
$ws.Range("E18").Value = '  -0.19%  '
$ws.Range("D19").Value = '''247.07'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.45%  '
$ws.Range("D20").Value = '''0.000007865'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.03%  '
$ws.Range("D21").Value = '''2.193.24'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.60%  '
$ws.Range("D22").Value = '''8.184'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.49%  '
$ws.Range("D23").Value = '''1.001'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("D25").Value = '''0.1595'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.60%  '
$ws.Range("D26").Value = '''9.539'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.41%  '
$ws.Range("D27").Value = '''162.24'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.26%  '
$ws.Range("D28").Value = '''18.81'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.14%  '
$ws.Range("D29").Value = '''2.044'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.92%  '
$ws.Range("E30").Value = '  +0.74%  '
$ws.Range("D31").Value = '''1.554'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.10%  '
$ws.Range("D32").Value = '''4.544'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.06%  '
$ws.Range("D33").Value = '''4.116'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.44%  '
$ws.Range("D34").Value = '''0.05485'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.89%  '
$ws.Range("D35").Value = '''1.249'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.61%  '
$ws.Range("D36").Value = '''0.7529'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.71%  '
$ws.Range("E37").Value = '  +0.25%  '
$ws.Range("D38").Value = '''2.690'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.04%  '
$ws.Range("D39").Value = '''0.01948'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.45%  '
$ws.Range("D40").Value = '''2.789'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.03%  '
$ws.Range("D41").Value = '''0.4503'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.33%  '
$ws.Range("D42").Value = '''74.38'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.95%  '
$ws.Range("D43").Value = '''6.089'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = '''1.094.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.19%  '
$ws.Range("D45").Value = '''0.8531'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.13%  '
$ws.Range("D46").Value = '''1.000'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("D47").Value = '''1.897'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.24%  '
$ws.Range("D48").Value = '''102.85'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.46%  '
$ws.Range("D49").Value = '''7.606'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.04%  '
$ws.Range("D50").Value = '''9.844'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.92%  '
$ws.Range("D51").Value = '''2.999'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.49%  '
